$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (shifts the existing "AT&T" row down to row 5)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new stock name
$ws.Range("A4").Value = "에이프로젠"

# Move selection to A6, matching the post-edit selection state
$ws.Range("A6").Select()
